# Update "想去人数" (column F) counts on the 展览 (sheet1), 演出 (sheet2)
# and 全部类型 (sheet4) worksheets to match the newly scraped totals.

$wb = $excel.ActiveWorkbook

# --- 展览 ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$sheet1Updates = @{
    2  = 10
    3  = 5127
    5  = 7430
    8  = 101
    11 = 26
    12 = 4311
    13 = 1757
    15 = 105
    16 = 2917
    19 = 205
    22 = 455
    24 = 98
    25 = 1690
    26 = 1184
    31 = 24
    34 = 60
    35 = 104
    37 = 2876
    38 = 702
    39 = 19
    40 = 64
    42 = 20
}
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Range("F$row").Value = $sheet1Updates[$row]
}

# --- 演出 -----------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$sheet2Updates = @{
    2 = 8
    3 = 10
}
foreach ($row in $sheet2Updates.Keys) {
    $ws2.Range("F$row").Value = $sheet2Updates[$row]
}

# --- 全部类型 -------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$sheet4Updates = @{
    2  = 10
    3  = 5127
    5  = 7430
    8  = 101
    11 = 26
    12 = 4311
    13 = 1757
    15 = 105
    16 = 2917
    19 = 205
    22 = 455
    23 = 8
    25 = 98
    26 = 1690
    27 = 1184
    32 = 24
    35 = 60
    36 = 104
    38 = 2876
    39 = 10
    40 = 702
    41 = 19
    42 = 64
    44 = 20
}
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Range("F$row").Value = $sheet4Updates[$row]
}
